$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the columns that are no longer part of the import template,
# from right-to-left so remaining column letters stay valid while we work.
# Original header layout (row 1):
#   A: userID                  <- remove
#   B: email
#   C: firstName
#   D: middleName
#   E: lastName
#   F: gender
#   G: rollNo
#   H: contactNo
#   I: alternateContactNo      <- remove
#   J: altenateContactPerson   <- remove
#   K: localAddress            <- remove
#   L: permanentAddress        <- remove
#   M: dob
#   N: linkedInProfile         <- remove
#   O: gitHubProfile           <- remove
#   P: blog                    <- remove

$ws.Range("P1").EntireColumn.Delete() | Out-Null
$ws.Range("O1").EntireColumn.Delete() | Out-Null
$ws.Range("N1").EntireColumn.Delete() | Out-Null
$ws.Range("L1").EntireColumn.Delete() | Out-Null
$ws.Range("K1").EntireColumn.Delete() | Out-Null
$ws.Range("J1").EntireColumn.Delete() | Out-Null
$ws.Range("I1").EntireColumn.Delete() | Out-Null
$ws.Range("A1").EntireColumn.Delete() | Out-Null

$ws.Range("G11").Select() | Out-Null
